$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that carries the "_GoBack" bookmark (currently the
# last, empty paragraph of the document). All new content is inserted
# immediately before and after it, exactly as the target diff shows.
# ---------------------------------------------------------------------------

function Get-GoBackParagraph {
    $bm = $d.Bookmarks("_GoBack")
    return $bm.Range.Paragraphs.First
}

# ---------------------------------------------------------------------------
# 1) Two new paragraphs BEFORE the bookmark paragraph (inserted in reading
#    order; each InsertParagraphBefore() lands right in front of _GoBack,
#    which is never bold, so the new paragraph marks stay non-bold too):
#      - "Présentation de l’application"   (bold heading)
#      - "Dans cette section ..."           (normal body text)
# ---------------------------------------------------------------------------

$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphBefore()
$gp = Get-GoBackParagraph
$p1 = $d.Paragraphs.Item($gp.Index - 1)
$p1.Range.InsertAfter("Présentation de l’application")
$p1.Range.Font.Bold = $true
$p1.Range.Font.BoldBi = $true

$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphBefore()
$gp = Get-GoBackParagraph
$p2 = $d.Paragraphs.Item($gp.Index - 1)
$p2.Range.InsertAfter("Dans cette section nous allons vous présenter, à l’aide de captures d’écrans, l’interface de l’application et justifier ses fonctionnalités.")

# ---------------------------------------------------------------------------
# 2) Five new paragraphs AFTER the bookmark paragraph:
#      - empty paragraph, bold paragraph mark
#      - "Annexes "                                    (bold heading)
#      - "Package windows : <dropbox url>"              (en-GB)
#      - "Github repository : <github url>"             (en-GB)
#
#    Each InsertParagraphAfter() call below is issued on the freshly
#    re-fetched _GoBack paragraph (never bold), so it always lands
#    immediately after _GoBack and the new paragraph never inherits bold
#    formatting. Building the block in REVERSE (last paragraph first)
#    therefore produces the correct final reading order without any bold
#    leaking from the "Annexes " heading into the two URL paragraphs.
# ---------------------------------------------------------------------------

# F) "Github repository : ..." (non-bold, en-GB)
$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphAfter()
$gp = Get-GoBackParagraph
$pF = $d.Paragraphs.Item($gp.Index + 1)
$pF.Range.InsertAfter("Github repository : ")
$pF.Range.InsertAfter("https://github.com/Ikewolf77/IHM_Projet/tree/main/Mini-Projet")
$pF.Range.LanguageID = "en-GB"

# E) "Package windows : ..." (non-bold, en-GB)
$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphAfter()
$gp = Get-GoBackParagraph
$pE = $d.Paragraphs.Item($gp.Index + 1)
$pE.Range.InsertAfter("Package windows : ")
$pE.Range.InsertAfter("https://www.dropbox.com/s/g0mj8mr8pxb6qm3/todo-list.zip?dl=0")
$pE.Range.LanguageID = "en-GB"

# D) "Annexes " (bold heading)
$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphAfter()
$gp = Get-GoBackParagraph
$pD = $d.Paragraphs.Item($gp.Index + 1)
$pD.Range.InsertAfter("Annexes ")
$pD.Range.Font.Bold = $true
$pD.Range.Font.BoldBi = $true

# C) empty paragraph, bold mark
$gp = Get-GoBackParagraph
$gp.Range.InsertParagraphAfter()
$gp = Get-GoBackParagraph
$pC = $d.Paragraphs.Item($gp.Index + 1)
$pC.Range.Font.Bold = $true
$pC.Range.Font.BoldBi = $true

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
